# lines_states.xlsx — add line7/line8 entries and shift the "extr" rows down,
# plus update several C/D/E values ("contingencies with rene fine").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Physically insert two new rows right after line6 (row 7), pushing the
# existing extr1..extr8 rows (old rows 8..15) down to rows 10..17, carrying
# their formatting/formulas with them.
$ws.Rows("8:9").Insert()

# Column A style (bold, centered, bordered) used by the index column, copied
# from an existing formatted cell so the two new rows match the rest.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A8:A9").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Target state for data rows 2..17 (row -> name, from_bus(C), to_bus(D), in_service(E))
$rows = [ordered]@{
    2  = @{ Name = "line1"; C = 7;  D = 9;  E = $true }
    3  = @{ Name = "line2"; C = 9;  D = 8;  E = $true }
    4  = @{ Name = "line3"; C = 8;  D = 10; E = $true }
    5  = @{ Name = "line4"; C = 8;  D = 11; E = $true }
    6  = @{ Name = "line5"; C = 10; D = 5;  E = $true }
    7  = @{ Name = "line6"; C = 12; D = 8;  E = $false }
    8  = @{ Name = "line7"; C = 14; D = 11; E = $true }
    9  = @{ Name = "line8"; C = 16; D = 9;  E = $true }
    10 = @{ Name = "extr1"; C = 5;  D = 12; E = $true }
    11 = @{ Name = "extr2"; C = 5;  D = 9;  E = $true }
    12 = @{ Name = "extr3"; C = 10; D = 11; E = $true }
    13 = @{ Name = "extr4"; C = 7;  D = 8;  E = $true }
    14 = @{ Name = "extr5"; C = 9;  D = 11; E = $true }
    15 = @{ Name = "extr6"; C = 7;  D = 11; E = $false }
    16 = @{ Name = "extr7"; C = 5;  D = 7;  E = $true }
    17 = @{ Name = "extr8"; C = 8;  D = 5;  E = $true }
}

foreach ($r in $rows.Keys) {
    $info = $rows[$r]

    # Column A: sequential index (row 2 -> 0, row 3 -> 1, ...)
    $ws.Cells.Item($r, 1).Value = ($r - 2)

    # Column B: name (shared string)
    $ws.Cells.Item($r, 2).Value = $info.Name

    # Columns C, D: numbers
    $ws.Cells.Item($r, 3).Value = $info.C
    $ws.Cells.Item($r, 4).Value = $info.D

    # Column E: boolean
    $ws.Cells.Item($r, 5).Value = $info.E
}
